# Update "new vacc rki" figures (Faktenblatt RKI Impf-Update 20.4. -> 21.4.)
$wb = $excel.ActiveWorkbook

# --- Sheet "Geimpfte Personen" ---------------------------------------------
$wsGeimpft = $wb.Worksheets.Item("Geimpfte Personen")
$wsGeimpft.Range("B1").Value = "Stand 21.4."
$wsGeimpft.Range("B3").Value = "17288804 (20,8 %)"
$wsGeimpft.Range("B4").Value = "11642016 (14,0 %)"
$wsGeimpft.Range("B5").Value = "5646788 ( 6,8 %)"

# --- Sheet "Regional Geimpfte" (Gesamt min. 1x / Gesamt 2x per Bundesland) --
$wsRegional = $wb.Worksheets.Item("Regional Geimpfte")

$wsRegional.Range("C2").Value = "20,8"
$wsRegional.Range("D2").Value = "6,8"

$wsRegional.Range("C3").Value = "19,9"
$wsRegional.Range("D3").Value = "6,7"

$wsRegional.Range("C4").Value = "21,4"
$wsRegional.Range("D4").Value = "6,8"

$wsRegional.Range("C5").Value = "20,2"
$wsRegional.Range("D5").Value = "8,1"

$wsRegional.Range("C6").Value = "20,2"
$wsRegional.Range("D6").Value = "6,7"

$wsRegional.Range("C7").Value = "23,5"
$wsRegional.Range("D7").Value = "7,4"

$wsRegional.Range("C8").Value = "21,0"
$wsRegional.Range("D8").Value = "5,8"

$wsRegional.Range("C9").Value = "18,8"
$wsRegional.Range("D9").Value = "7,3"

$wsRegional.Range("C10").Value = "20,8"
$wsRegional.Range("D10").Value = "6,4"

$wsRegional.Range("C11").Value = "20,0"
$wsRegional.Range("D11").Value = "6,4"

$wsRegional.Range("C12").Value = "21,6"
$wsRegional.Range("D12").Value = "6,6"

$wsRegional.Range("C13").Value = "21,0"
$wsRegional.Range("D13").Value = "6,4"

$wsRegional.Range("C14").Value = "23,2"
$wsRegional.Range("D14").Value = "7,7"

$wsRegional.Range("C15").Value = "20,4"
$wsRegional.Range("D15").Value = "7,4"

$wsRegional.Range("C16").Value = "21,6"
$wsRegional.Range("D16").Value = "5,8"

$wsRegional.Range("C17").Value = "21,3"
$wsRegional.Range("D17").Value = "6,2"

$wsRegional.Range("C18").Value = "20,6"
$wsRegional.Range("D18").Value = "8,3"

# --- Sheet "Impfstoffdosen" (dose counts per manufacturer) ------------------
$wsDosen = $wb.Worksheets.Item("Impfstoffdosen")

$wsDosen.Range("B2").Value = "16905070 ( 96,0 %)"
$wsDosen.Range("B3").Value = "11571347"
$wsDosen.Range("B4").Value = "5333723"

$wsDosen.Range("B6").Value = "1242388 (  7,1 %)"
$wsDosen.Range("B7").Value = "939534"
$wsDosen.Range("B8").Value = "302854"

$wsDosen.Range("B10").Value = "4788134 ( 27,2 %)"
$wsDosen.Range("B11").Value = "4777923"
$wsDosen.Range("B12").Value = "10211"
